# Applies the 2016-11-16 protocol edits:
#   1. "Selina Brinnich" participant name gets spell-check markup
#      (w:proofErr spellStart/spellEnd) wrapped around "Brinnich",
#      splitting the surrounding run into three runs.
#   2. The "next meeting" date/time moves from 15.11.2016 (um 10:30)
#      to 30.11.2016 (um 9:00), with the run layout changing and the
#      "_GoBack" bookmark shifting from after "FH" to inside "(um 9:00)".

$d = $word.ActiveDocument

# --- Package template used for the surgical InsertXML calls below ---
# InsertXML() on a Range deletes that range's content and appends the
# supplied run/field-level markup to the end of the paragraph that
# contained the (now collapsed) range, so each replacement range below
# is deliberately extended out to the paragraph's existing end so the
# untouched tail text ("A 5.35", etc.) can be re-supplied in the new
# markup and land back in the correct order.
$pkgOpen = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$pkgClose = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Edit 1: Teilnehmer list -- wrap "Brinnich" with spell-check proofErr
# ---------------------------------------------------------------------
# NOTE: Find.Execute's Range must be re-materialized via $d.Range(start,end)
# before InsertXML -- calling InsertXML directly on the Range returned/
# mutated by Find leaves the matched text in place instead of replacing it.
$namesFindRange = $d.Content
$namesFound = $namesFindRange.Find.Execute(", Phillip Schermann, Selina Brinnich, Thiago Gumhold", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($namesFound) {
    $namesRange = $d.Range($namesFindRange.Start, $namesFindRange.End)
    $namesXml = $pkgOpen + `
        '<w:r><w:t xml:space="preserve">, Phillip Schermann, Selina </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>Brinnich</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t>, Thiago Gumhold</w:t></w:r>' + `
        $pkgClose
    $namesRange.InsertXML($namesXml)
}

# ---------------------------------------------------------------------
# Edit 2: Next meeting date/time + relocated "_GoBack" bookmark
# ---------------------------------------------------------------------
$meetingAnchor = $d.Content
$meetingFound = $meetingAnchor.Find.Execute("Meeting am", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterMeetingAm = $meetingAnchor.End

$tailAnchor = $d.Range($afterMeetingAm, $d.Content.End)
$tailFound = $tailAnchor.Find.Execute("A 5.35", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endOfA535 = $tailAnchor.End

if ($meetingFound -and $tailFound) {
    $dateRange = $d.Range($afterMeetingAm, $endOfA535)

    $enDash = [char]0x2013

    $dateXml = $pkgOpen + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:r><w:t>30</w:t></w:r>' + `
        '<w:r><w:t>.11.2016</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> (um 9:0</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
        '<w:bookmarkEnd w:id="0"/>' + `
        '<w:r><w:t>0)</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:r><w:t>' + $enDash + '</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> FH</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:r><w:t>A 5.35</w:t></w:r>' + `
        $pkgClose

    $dateRange.InsertXML($dateXml)
}
